$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated figures for existing rows (columns C, D, I, N, O, P) ---
$ws.Range("O198").Value = 212
$ws.Range("O199").Value = 226
$ws.Range("O201").Value = 257
$ws.Range("O202").Value = 295
$ws.Range("O203").Value = 289
$ws.Range("O204").Value = 331
$ws.Range("O205").Value = 385
$ws.Range("O206").Value = 408
$ws.Range("O207").Value = 427
$ws.Range("O208").Value = 420
$ws.Range("O209").Value = 385

$ws.Range("I210").Value = 1
$ws.Range("N210").Value = 118
$ws.Range("O210").Value = 399

$ws.Range("N211").Value = 116

$ws.Range("C212").Value = 11

$ws.Range("O213").Value = 446

$ws.Range("O214").Value = 392
$ws.Range("P214").Value = 518

$ws.Range("O215").Value = 349
$ws.Range("P215").Value = 532

$ws.Range("C216").Value = 21
$ws.Range("N216").Value = 81
$ws.Range("O216").Value = 309
$ws.Range("P216").Value = 565

$ws.Range("C217").Value = 17
$ws.Range("D217").Value = 2
$ws.Range("N217").Value = 98
$ws.Range("O217").Value = 300
$ws.Range("P217").Value = 548

$ws.Range("C218").Value = 9
$ws.Range("N218").Value = 96
$ws.Range("O218").Value = 251
$ws.Range("P218").Value = 480

# --- Fill in the newly-added data row 219 (2020-10-01) ---
$ws.Range("C219").Value = 0
$ws.Range("D219").Value = 0
$ws.Range("E219").Value = 1
$ws.Range("F219").Value = 1
$ws.Range("G219").Value = 7
$ws.Range("I219").Value = 0
$ws.Range("L219").Value = "0"
$ws.Range("M219").Value = "0"
$ws.Range("N219").Value = 91
$ws.Range("O219").Value = 214
$ws.Range("P219").Value = 443

# Recalculate the running-total formulas (columns B, H, J, K)
$excel.Calculate()

# --- View state: frozen pane top-left cell and active selection ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("P229").Select()
